$wb = $excel.ActiveWorkbook

# This script applies a scheduled market-data refresh to the Tonberry Profits
# workbook: updated currentAveragePrice / NQ / HQ price columns (H, I, J, K, L)
# and the recomputed profit columns (M = LeveProfitNQ, N = LeveProfitHQ) for the
# rows whose underlying market data changed. A few rows lose their NQ or HQ price
# entirely (price data no longer available), so the corresponding profit cell is
# cleared rather than set.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 11).Value = 91
$ws.Cells.Item(33, 13).Value = 138
$ws.Cells.Item(33, 9).Value = 91
$ws.Cells.Item(33, 8).Value = 107
$ws.Cells.Item(95, 12).Value = 37999.5
$ws.Cells.Item(95, 10).Value = 37999.5
$ws.Cells.Item(95, 14).Value = -43491.5
$ws.Cells.Item(95, 8).Value = 37999.5
$ws.Cells.Item(107, 12).Value = 2000.4
$ws.Cells.Item(107, 11).Value = 591.7273
$ws.Cells.Item(107, 13).Value = 1328.2727
$ws.Cells.Item(107, 10).Value = 2000.4
$ws.Cells.Item(107, 9).Value = 591.7273
$ws.Cells.Item(107, 14).Value = -5840.4
$ws.Cells.Item(107, 8).Value = 1031.9375
$ws.Cells.Item(116, 12).Value = 3133.3333
$ws.Cells.Item(116, 11).Value = 22256.4
$ws.Cells.Item(116, 13).Value = -18814.4
$ws.Cells.Item(116, 10).Value = 3133.3333
$ws.Cells.Item(116, 9).Value = 22256.4
$ws.Cells.Item(116, 14).Value = -10017.3333
$ws.Cells.Item(116, 8).Value = 9963
$ws.Cells.Item(132, 11).Value = 3103.3125
$ws.Cells.Item(132, 13).Value = -573.3125
$ws.Cells.Item(132, 9).Value = 1034.4375
$ws.Cells.Item(132, 8).Value = 1034.4375
$ws.Cells.Item(137, 11).Value = 5977.5
$ws.Cells.Item(137, 13).Value = -3427.5
$ws.Cells.Item(137, 9).Value = 1992.5
$ws.Cells.Item(137, 8).Value = 2135
$ws.Cells.Item(139, 12).Value = 48100
$ws.Cells.Item(139, 10).Value = 48100
$ws.Cells.Item(139, 14).Value = -58380
$ws.Cells.Item(139, 8).Value = 48100
$ws.Cells.Item(140, 12).Value = 53335
$ws.Cells.Item(140, 10).Value = 53335
$ws.Cells.Item(140, 14).Value = -63695
$ws.Cells.Item(140, 8).Value = 53335

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 12).Value = 6307.9287
$ws.Cells.Item(32, 10).Value = 6307.9287
$ws.Cells.Item(32, 14).Value = -6881.9287
$ws.Cells.Item(32, 8).Value = 3154.2188
$ws.Cells.Item(61, 11).Value = 1700.7778
$ws.Cells.Item(61, 13).Value = -1488.7778
$ws.Cells.Item(61, 9).Value = 1700.7778
$ws.Cells.Item(61, 8).Value = 2403.7
$ws.Cells.Item(74, 11).Value = 1480.7
$ws.Cells.Item(74, 13).Value = -606.7
$ws.Cells.Item(74, 9).Value = 1480.7
$ws.Cells.Item(74, 8).Value = 1572.2727
$ws.Cells.Item(77, 11).Value = 7403.5
$ws.Cells.Item(77, 13).Value = -3035.5
$ws.Cells.Item(77, 9).Value = 1480.7
$ws.Cells.Item(77, 8).Value = 1572.2727
$ws.Cells.Item(122, 11).Value = 5103
$ws.Cells.Item(122, 13).Value = -2653
$ws.Cells.Item(122, 9).Value = 1701
$ws.Cells.Item(122, 8).Value = 1679.2174
$ws.Cells.Item(132, 12).Value = 11998.2
$ws.Cells.Item(132, 11).Value = 10978.2
$ws.Cells.Item(132, 13).Value = -8448.200000000001
$ws.Cells.Item(132, 10).Value = 3999.4
$ws.Cells.Item(132, 9).Value = 3659.4
$ws.Cells.Item(132, 14).Value = -17058.2
$ws.Cells.Item(132, 8).Value = 3829.4
$ws.Cells.Item(136, 11).Value = 5102.3334
$ws.Cells.Item(136, 13).Value = -2552.3334
$ws.Cells.Item(136, 9).Value = 1700.7778
$ws.Cells.Item(136, 8).Value = 2403.7

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(105, 11).Value = 1902.0952
$ws.Cells.Item(105, 13).Value = -155.0952
$ws.Cells.Item(105, 9).Value = 1902.0952
$ws.Cells.Item(105, 8).Value = 1905.3462
$ws.Cells.Item(28, 14).ClearContents()
$ws.Cells.Item(98, 14).ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 12).Value = 45599
$ws.Cells.Item(20, 10).Value = 45599
$ws.Cells.Item(20, 14).Value = -46071
$ws.Cells.Item(20, 8).Value = 45599
$ws.Cells.Item(30, 12).Value = 45599
$ws.Cells.Item(30, 10).Value = 45599
$ws.Cells.Item(30, 14).Value = -45781
$ws.Cells.Item(30, 8).Value = 45599
$ws.Cells.Item(31, 12).Value = 7689
$ws.Cells.Item(31, 11).Value = 1769.3334
$ws.Cells.Item(31, 13).Value = -1474.3334
$ws.Cells.Item(31, 10).Value = 7689
$ws.Cells.Item(31, 9).Value = 1769.3334
$ws.Cells.Item(31, 14).Value = -8279
$ws.Cells.Item(31, 8).Value = 5152
$ws.Cells.Item(34, 12).Value = 7689
$ws.Cells.Item(34, 11).Value = 1769.3334
$ws.Cells.Item(34, 13).Value = -1567.3334
$ws.Cells.Item(34, 10).Value = 7689
$ws.Cells.Item(34, 9).Value = 1769.3334
$ws.Cells.Item(34, 14).Value = -8093
$ws.Cells.Item(34, 8).Value = 5152
$ws.Cells.Item(58, 11).Value = 1142.2
$ws.Cells.Item(58, 13).Value = -939.2
$ws.Cells.Item(58, 9).Value = 1142.2
$ws.Cells.Item(58, 8).Value = 1101.4286
$ws.Cells.Item(59, 12).Value = 39450
$ws.Cells.Item(59, 10).Value = 39450
$ws.Cells.Item(59, 14).Value = -41740
$ws.Cells.Item(59, 8).Value = 39450
$ws.Cells.Item(62, 11).Value = 3799.4
$ws.Cells.Item(62, 13).Value = -3175.4
$ws.Cells.Item(62, 9).Value = 3799.4
$ws.Cells.Item(62, 8).Value = 3785.2856
$ws.Cells.Item(65, 11).Value = 18997
$ws.Cells.Item(65, 13).Value = -15877
$ws.Cells.Item(65, 9).Value = 3799.4
$ws.Cells.Item(65, 8).Value = 3785.2856
$ws.Cells.Item(94, 11).Value = 1201.7142
$ws.Cells.Item(94, 13).Value = -750.7141999999999
$ws.Cells.Item(94, 9).Value = 1201.7142
$ws.Cells.Item(94, 8).Value = 1161.2
$ws.Cells.Item(122, 12).Value = 8467.200000000001
$ws.Cells.Item(122, 10).Value = 2822.4
$ws.Cells.Item(122, 14).Value = -13367.2
$ws.Cells.Item(122, 8).Value = 2811.4
$ws.Cells.Item(128, 12).Value = 45599
$ws.Cells.Item(128, 10).Value = 45599
$ws.Cells.Item(128, 14).Value = -55559
$ws.Cells.Item(128, 8).Value = 45599
$ws.Cells.Item(136, 11).Value = 3426.6
$ws.Cells.Item(136, 13).Value = -876.6000000000004
$ws.Cells.Item(136, 9).Value = 1142.2
$ws.Cells.Item(136, 8).Value = 1101.4286

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 12).Value = 5259.4998
$ws.Cells.Item(50, 10).Value = 1753.1666
$ws.Cells.Item(50, 14).Value = -6221.4998
$ws.Cells.Item(50, 8).Value = 78765.89
$ws.Cells.Item(53, 12).Value = 5259.4998
$ws.Cells.Item(53, 10).Value = 1753.1666
$ws.Cells.Item(53, 14).Value = -6221.4998
$ws.Cells.Item(53, 8).Value = 78765.89
$ws.Cells.Item(130, 12).Value = 6600
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 10).Value = 2200
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 14).Value = -16640
$ws.Cells.Item(130, 8).Value = 2200
$ws.Cells.Item(131, 12).Value = 300886.29
$ws.Cells.Item(131, 11).Value = 150001500
$ws.Cells.Item(131, 13).Value = -149996460
$ws.Cells.Item(131, 10).Value = 100295.43
$ws.Cells.Item(131, 9).Value = 50000500
$ws.Cells.Item(131, 14).Value = -310966.29
$ws.Cells.Item(131, 8).Value = 29453358
$ws.Cells.Item(137, 11).Value = 3822.6666
$ws.Cells.Item(137, 13).Value = 1277.3334
$ws.Cells.Item(137, 9).Value = 1274.2222
$ws.Cells.Item(137, 8).Value = 3538.1538
$ws.Cells.Item(130, 13).ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(104, 12).Value = 49850
$ws.Cells.Item(104, 10).Value = 49850
$ws.Cells.Item(104, 14).Value = -56838
$ws.Cells.Item(104, 8).Value = 49850
$ws.Cells.Item(110, 12).Value = 99792
$ws.Cells.Item(110, 10).Value = 99792
$ws.Cells.Item(110, 14).Value = -107972
$ws.Cells.Item(110, 8).Value = 99792
$ws.Cells.Item(134, 12).Value = 130711.71
$ws.Cells.Item(134, 10).Value = 43570.57
$ws.Cells.Item(134, 14).Value = -135781.71
$ws.Cells.Item(134, 8).Value = 43570.57

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 11).Value = 1063.5714
$ws.Cells.Item(93, 13).Value = 184.4286
$ws.Cells.Item(93, 9).Value = 1063.5714
$ws.Cells.Item(93, 8).Value = 1269.4445

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 12).Value = 47452.75
$ws.Cells.Item(123, 10).Value = 47452.75
$ws.Cells.Item(123, 14).Value = -57252.75
$ws.Cells.Item(123, 8).Value = 47452.75

Write-Output "Applied Tonberry_Profits market-data refresh"